$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old table (A1:D4) - both values and formatting - first, since
# the sheet is being rebuilt with a new layout (A1:E7): the "code" column
# moves from A to C, "genderName"->name moves to D, "langCode"->lang_code
# moves to B, "isActive"->is_active moves to E, and a new numeric id column
# is introduced in A.
$ws.Range("A1:D4").Clear() | Out-Null

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "lang_code"
$ws.Range("C1").Value = "code"
$ws.Range("D1").Value = "name"
$ws.Range("E1").Value = "is_active"

# ---- Data rows (2-7): id, lang_code, code, name, is_active ----
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "eng"
$ws.Range("C2").Value = "MLE"
$ws.Range("D2").Value = "Male"
$ws.Range("E2").Value = $true

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "eng"
$ws.Range("C3").Value = "FLE"
$ws.Range("D3").Value = "Female"
$ws.Range("E3").Value = $true

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "eng"
$ws.Range("C4").Value = "OTH"
$ws.Range("D4").Value = "Others"
$ws.Range("E4").Value = $false

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "fra"
$ws.Range("C5").Value = "MLE"
$ws.Range("D5").Value = "Mâle"
$ws.Range("E5").Value = $true

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "fra"
$ws.Range("C6").Value = "FLE"
$ws.Range("D6").Value = "Femelle"
$ws.Range("E6").Value = $true

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "fra"
$ws.Range("C7").Value = "OTH"
$ws.Range("D7").Value = "Dautres"
$ws.Range("E7").Value = $false

# ---- Re-apply the bold / bordered / centered header style (originally
# style index 1, used only by A1 "code") across the whole new header row
# (B1:E1) and the id column (A2:A7), which keep that same look. Applied as
# two single-area ranges since multi-area Range property writes only hit
# the first area in this host.
$headerRow = $ws.Range("B1:E1")
$headerRow.Font.Bold = $true
$headerRow.HorizontalAlignment = -4108   # xlCenter
$headerRow.VerticalAlignment = -4160     # xlTop
$headerRow.Borders.LineStyle = 1         # xlContinuous
$headerRow.Borders.Weight = 2            # xlThin

$idCol = $ws.Range("A2:A7")
$idCol.Font.Bold = $true
$idCol.HorizontalAlignment = -4108       # xlCenter
$idCol.VerticalAlignment = -4160         # xlTop
$idCol.Borders.LineStyle = 1             # xlContinuous
$idCol.Borders.Weight = 2                # xlThin
